# Apply updated enrollment counts (Inscritos/Pagos/Inscrições homologadas)
# to the "Inscricoes" sheet, per the latest registration report figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Edificacoes - Campus Congonhas): Pagos / Inscrições homologadas 10 -> 11
$ws.Range("F3").Value = 11
$ws.Range("H3").Value = 11

# Row 6 (Eletrotecnica - Campus Conselheiro Lafaiete): 25 -> 26
$ws.Range("F6").Value = 26
$ws.Range("H6").Value = 26

# Row 7 (Mecanica - Campus Conselheiro Lafaiete): 14 -> 15
$ws.Range("F7").Value = 15
$ws.Range("H7").Value = 15

# Row 10: 12 -> 13
$ws.Range("F10").Value = 13
$ws.Range("H10").Value = 13

# Row 15: Inscritos 88 -> 89, Pagos / Inscrições homologadas 41 -> 43
$ws.Range("E15").Value = 89
$ws.Range("F15").Value = 43
$ws.Range("H15").Value = 43

# Row 16: Inscritos 300 -> 301 (Pagos / homologadas unchanged)
$ws.Range("E16").Value = 301
